$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price-column cells whose new value would otherwise be
# auto-converted to a number by Excel (values with a single "." look numeric).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated Price / Volume(1h) values scraped by the cron job.
$ws.Range("D2").Value = '26.521.64'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '1.627.69'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '213.05'
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D10").Value = '18.79'
$ws.Range("E10").Value = '  -1.32%  '
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("D12").Value = '1.852.82'
$ws.Range("E12").Value = '  -0.54%  '
$ws.Range("D13").Value = '1.634.18'
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '64.98'
$ws.Range("E16").Value = '  +3.09%  '
$ws.Range("D17").Value = '26.526.98'
$ws.Range("E17").Value = '  -0.56%  '
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = '214.38'
$ws.Range("E19").Value = '  +2.74%  '
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("D22").Value = '6.25'
$ws.Range("E22").Value = '  +1.42%  '
$ws.Range("D23").Value = '9.29'
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("E24").Value = '  +9.26%  '
$ws.Range("D25").Value = '148.53'
$ws.Range("E25").Value = '  +1.43%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  -0.29%  '
$ws.Range("D28").Value = '6.87'
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("D29").Value = '15.53'
$ws.Range("E29").Value = '  +0.93%  '
$ws.Range("E30").Value = '  -1.74%  '
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("E32").Value = '  +3.21%  '
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("D34").Value = '1.238.37'
$ws.Range("E34").Value = '  +5.93%  '
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("E37").Value = '  +4.18%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = '0.507'
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("D40").Value = '0.795'
$ws.Range("E40").Value = '  -1.49%  '
$ws.Range("E41").Value = '  -1.94%  '
$ws.Range("D42").Value = '0.799'
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("D44").Value = '1.763.60'
$ws.Range("E44").Value = '  -0.58%  '
$ws.Range("D45").Value = '92.92'
$ws.Range("E45").Value = '  +0.61%  '
$ws.Range("E46").Value = '  +2.55%  '
$ws.Range("D47").Value = '54.92'
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("E48").Value = '  -0.72%  '
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("D50").Value = '7.48'
$ws.Range("E50").Value = '  -0.80%  '
$ws.Range("E51").Value = '  +0.31%  '
